# This script updates the "cryptos" worksheet (Coin / Link / Price / Volume(1h))
# to reflect the latest scraped values from coinranking.com, per the
# "Updated cryptos list ... with GitHub Actions" commit.
#
# For each affected row we update Coin (B), Link (C), Price (D) and/or
# Volume(1h) (E) as needed. Three row-pairs (12/13, 37/38, 41/42) have their
# Coin/Link/Price/Volume swapped because the two coins changed relative rank.
#
# Price values are textual (e.g. "28.275.15", "1.000") so we force the
# NumberFormat to Text ("@") before assigning them; otherwise Excel would
# auto-convert them to numbers and silently drop formatting (trailing
# zeros, multiple "thousand separator" dots, etc.).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2; B=$null; C=$null; D='28.275.15'; E='  +1.06%  '},
    @{Row=3; B=$null; C=$null; D='1.804.41'; E='  +2.89%  '},
    @{Row=4; B=$null; C=$null; D='1.002'; E='  -0.09%  '},
    @{Row=5; B=$null; C=$null; D='338.24'; E='  +0.75%  '},
    @{Row=6; B=$null; C=$null; D='1.000'; E='  +0.12%  '},
    @{Row=7; B=$null; C=$null; D='0.4678'; E='  +21.94%  '},
    @{Row=8; B=$null; C=$null; D='0.3794'; E='  +10.88%  '},
    @{Row=9; B=$null; C=$null; D='45.51'; E='  -0.69%  '},
    @{Row=10; B=$null; C=$null; D='1.150'; E='  +2.61%  '},
    @{Row=11; B=$null; C=$null; D='0.07609'; E='  +5.09%  '},
    @{Row=12; B='BinanceUSD'; C='https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'; D='1.001'; E='  -0.04%  '},
    @{Row=13; B='Solana'; C='https://coinranking.com/coin/zNZHO_Sjf+solana-sol'; D='22.36'; E='  -0.82%  '},
    @{Row=14; B=$null; C=$null; D='6.324'; E='  +2.45%  '},
    @{Row=15; B=$null; C=$null; D='7.447'; E='  +4.11%  '},
    @{Row=16; B=$null; C=$null; D='1.805.33'; E='  +3.16%  '},
    @{Row=17; B=$null; C=$null; D='0.00001093'; E='  +2.94%  '},
    @{Row=18; B=$null; C=$null; D='0.06728'; E='  +1.90%  '},
    @{Row=19; B=$null; C=$null; D='81.74'; E='  +3.04%  '},
    @{Row=20; B=$null; C=$null; D='0.9993'; E='  +0.01%  '},
    @{Row=21; B=$null; C=$null; D='17.41'; E='  +3.88%  '},
    @{Row=22; B=$null; C=$null; D='6.404'; E='  +3.30%  '},
    @{Row=23; B=$null; C=$null; D='28.253.76'; E='  +0.94%  '},
    @{Row=24; B=$null; C=$null; D='11.84'; E='  +1.36%  '},
    @{Row=25; B=$null; C=$null; D='2.408'; E='  +1.37%  '},
    @{Row=26; B=$null; C=$null; D='20.72'; E='  +4.27%  '},
    @{Row=27; B=$null; C=$null; D='153.91'; E='  -0.28%  '},
    @{Row=28; B=$null; C=$null; D=$null; E='  +3.05%  '},
    @{Row=29; B=$null; C=$null; D='2.009.82'; E='  +3.05%  '},
    @{Row=30; B=$null; C=$null; D='133.06'; E='  +1.35%  '},
    @{Row=31; B=$null; C=$null; D='1.251'; E='  -0.49%  '},
    @{Row=32; B=$null; C=$null; D='4.037'; E='  +0.21%  '},
    @{Row=33; B=$null; C=$null; D='0.09624'; E='  +9.17%  '},
    @{Row=34; B=$null; C=$null; D='5.851'; E='  -0.02%  '},
    @{Row=35; B=$null; C=$null; D='0.2235'; E='  +6.15%  '},
    @{Row=36; B=$null; C=$null; D='0.06370'; E='  +3.48%  '},
    @{Row=37; B='Aptos'; C='https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; D='12.10'; E='  -1.01%  '},
    @{Row=38; B='VeChain'; C='https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; D='0.02351'; E='  +2.93%  '},
    @{Row=39; B=$null; C=$null; D='5.250'; E='  +2.06%  '},
    @{Row=40; B=$null; C=$null; D='0.6621'; E='  +0.73%  '},
    @{Row=41; B='TrustWalletToken'; C='https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; D='1.237'; E='  +1.78%  '},
    @{Row=42; B='WEMIXTOKEN'; C='https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'; D='1.502'; E='  -2.92%  '},
    @{Row=43; B=$null; C=$null; D=$null; E='  +3.45%  '},
    @{Row=44; B=$null; C=$null; D='14.17'; E='  +3.22%  '},
    @{Row=45; B=$null; C=$null; D='0.9998'; E=$null},
    @{Row=46; B=$null; C=$null; D='0.6119'; E='  +1.07%  '},
    @{Row=47; B=$null; C=$null; D='3.844'; E='  +0.15%  '},
    @{Row=48; B=$null; C=$null; D='130.51'; E='  +2.57%  '},
    @{Row=49; B=$null; C=$null; D='2.033'; E='  +1.28%  '},
    @{Row=50; B=$null; C=$null; D='0.07159'; E='  +2.58%  '},
    @{Row=51; B=$null; C=$null; D=$null; E='  +0.62%  '}

)

foreach ($u in $updates) {
    $r = $u.Row

    if ($null -ne $u.B) {
        $ws.Cells.Item($r, 2).Value = $u.B
    }
    if ($null -ne $u.C) {
        $ws.Cells.Item($r, 3).Value = $u.C
    }
    if ($null -ne $u.D) {
        $dCell = $ws.Cells.Item($r, 4)
        $dCell.NumberFormat = "@"
        $dCell.Value = $u.D
    }
    if ($null -ne $u.E) {
        $ws.Cells.Item($r, 5).Value = $u.E
    }
}
